# Add a new "2023" column (column S) to the hotels and restaurants sheet,
# mirroring the formatting of the existing "2022" column (column R) for
# each row, and update the selected cell in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> value for the new 2023 column (S)
$values = @{
    3  = 2023
    4  = 130.9
    5  = 131.1
    6  = 2047
    7  = 2003
    8  = 1496.7
    9  = 78.1
    10 = 36.5
    11 = 53
    12 = 24.2
    13 = 77.2
    14 = 0.1567
}

foreach ($row in 3..14) {
    $src = $ws.Cells.Item($row, 18)   # column R
    $dst = $ws.Cells.Item($row, 19)   # column S

    # Copy the source cell's formatting (number format, font, borders,
    # alignment, etc.) onto the new cell before writing its value.
    $src.Copy()
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    $dst.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Update the active selection to match the saved workbook state.
[void]$ws.Range("F19").Select()
